# Apply "Storefront updates, File Upload section added, Versions section added"
# to the docs_virtocommerce_org-sitemap worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CLI-tools section (rows 13-17): new target URLs + status moved to OnReview ---
$ws.Range("B13").Value = "platform/developer-guide/CLI-tools/build-automation/"
$ws.Range("C13").Value = "OnReview"

$ws.Range("B14").Value = "platform/developer-guide/CLI-tools/cold-start-data-migration/"
$ws.Range("C14").Value = "OnReview"

$ws.Range("B15").Value = "platform/developer-guide/CLI-tools/grab-migrator/"
$ws.Range("C15").Value = "OnReview"

$ws.Range("B16").Value = "platform/developer-guide/CLI-tools/overview/"
$ws.Range("C16").Value = "OnReview"

$ws.Range("B17").Value = "platform/developer-guide/CLI-tools/package-management/"
$ws.Range("C17").Value = "OnReview"

# --- Versions section (rows 169-180): new target URLs, status OnReview, responsible Maria ---
$ws.Range("B169").Value = "platform/user-guide/versions/virto3-products-versions/"
$ws.Range("C169").Value = "OnReview"
$ws.Range("D169").Value = "Maria"

$ws.Range("B170").Value = "platform/user-guide/versions/v3-2020-Q3/"
$ws.Range("C170").Value = "OnReview"
$ws.Range("D170").Value = "Maria"

$ws.Range("B171").Value = "platform/user-guide/versions/v3-2021-Q1/"
$ws.Range("C171").Value = "OnReview"
$ws.Range("D171").Value = "Maria"

$ws.Range("B172").Value = "platform/user-guide/versions/v3-2021-Q3/"
$ws.Range("C172").Value = "OnReview"
$ws.Range("D172").Value = "Maria"

$ws.Range("B173").Value = "platform/user-guide/versions/v3-2022-Q1/"
$ws.Range("C173").Value = "OnReview"
$ws.Range("D173").Value = "Maria"

$ws.Range("B174").Value = "platform/user-guide/versions/v3-2022-Q2/"
$ws.Range("C174").Value = "OnReview"
$ws.Range("D174").Value = "Maria"

$ws.Range("B175").Value = "platform/user-guide/versions/v3-2022-Q3/"
$ws.Range("C175").Value = "OnReview"
$ws.Range("D175").Value = "Maria"

$ws.Range("B176").Value = "platform/user-guide/versions/v3-2022-Q4/"
$ws.Range("C176").Value = "OnReview"
$ws.Range("D176").Value = "Maria"

$ws.Range("B177").Value = "platform/user-guide/versions/v3-2023-S5/"
$ws.Range("C177").Value = "OnReview"
$ws.Range("D177").Value = "Maria"

$ws.Range("B178").Value = "platform/user-guide/versions/v3-2023-S6/"
$ws.Range("C178").Value = "OnReview"
$ws.Range("D178").Value = "Maria"
$ws.Range("B178").Style = "Hyperlink"

$ws.Range("B179").Value = "platform/user-guide/versions/v3-2024-S7/"
$ws.Range("C179").Value = "OnReview"
$ws.Range("D179").Value = "Maria"

$ws.Range("B180").Value = "platform/user-guide/versions/v3-2024-S8/"
$ws.Range("C180").Value = "OnReview"
$ws.Range("D180").Value = "Maria"

# --- restore cursor/selection position as last saved by the author ---
$ws.Range("B181").Select()
$excel.ActiveWindow.ScrollRow = 154
$excel.ActiveWindow.ScrollColumn = 1
